$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("NewsagentUserStories")
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 27
$ws2.Range("A1:F48").Select()
Write-Host "ScrollRow final:" $excel.ActiveWindow.ScrollRow
